$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3306367468558449
$ws.Range("C2").Value = 0.06440320679725176
$ws.Range("D2").Value = 0.02380209259329291
$ws.Range("E2").Value = 0.4180825186005563
$ws.Range("F2").Value = 0.6122752868145653
$ws.Range("I2").Value = 0.4517866951497496
$ws.Range("K2").Value = 0.3611763121003264
$ws.Range("N2").Value = 1.189202072763843
$ws.Range("O2").Value = 2.053608104493634
$ws.Range("B3").Value = 0.2911160902318386
$ws.Range("C3").Value = 0.05667629207731295
$ws.Range("D3").Value = 0.02184017161742702
$ws.Range("E3").Value = 0.3648139447740704
$ws.Range("F3").Value = 0.6089813308479037
$ws.Range("I3").Value = 0.4552154285460013
$ws.Range("K3").Value = 0.315793928048123
$ws.Range("N3").Value = 1.203388477486316
$ws.Range("O3").Value = 2.05573438518465
$ws.Range("B4").Value = 0.2668413665212199
$ws.Range("C4").Value = 0.05190885488340768
$ws.Range("D4").Value = 0.0206255991230222
$ws.Range("E4").Value = 0.3321876317310313
$ws.Range("F4").Value = 0.6073720781821166
$ws.Range("I4").Value = 0.4575912503754189
$ws.Range("K4").Value = 0.2878860502833902
$ws.Range("N4").Value = 1.212544879463678
$ws.Range("O4").Value = 2.058376677175445
$ws.Range("B5").Value = 0.2569474271387833
$ws.Range("C5").Value = 0.04996030176565114
$ws.Range("D5").Value = 0.02012818275431272
$ws.Range("E5").Value = 0.3189109501738301
$ws.Range("F5").Value = 0.6068201013650665
$ws.Range("I5").Value = 0.458627377340779
$ws.Range("K5").Value = 0.2765029866158386
$ws.Range("N5").Value = 1.216388262846045
$ws.Range("O5").Value = 2.059789110565589
$ws.Range("B6").Value = 0.2553044507992297
$ws.Range("C6").Value = 0.04963639728438807
$ws.Range("D6").Value = 0.0200454390299214
$ws.Range("E6").Value = 0.3167074590422914
$ws.Range("F6").Value = 0.6067347125572269
$ws.Range("I6").Value = 0.4588035278624396
$ws.Range("K6").Value = 0.2746122230464607
$ws.Range("N6").Value = 1.217033221440543
$ws.Range("O6").Value = 2.060043904841578
$ws.Range("B7").Value = 0.2667079400669081
$ws.Range("C7").Value = 0.05188259940449313
$ws.Range("D7").Value = 0.02061890074220685
$ws.Range("E7").Value = 0.3320085036035749
$ws.Range("F7").Value = 0.6073642138654947
$ws.Range("I7").Value = 0.4576049488790694
$ws.Range("K7").Value = 0.2877325757725941
$ws.Range("N7").Value = 1.212596258953075
$ws.Range("O7").Value = 2.05839436721304
$ws.Range("B8").Value = 0.3170121197505864
$ws.Range("C8").Value = 0.06174376659831182
$ws.Range("D8").Value = 0.0231277052974832
$ws.Range("E8").Value = 0.3996978156532123
$ws.Range("F8").Value = 0.6110536926471255
$ws.Range("I8").Value = 0.4529127206142718
$ws.Range("K8").Value = 0.3455376061479001
$ws.Range("N8").Value = 1.194000972432589
$ws.Range("O8").Value = 2.054063522378527
$ws.Range("B9").Value = 0.4155734735244891
$ws.Range("C9").Value = 0.08089839569993273
$ws.Range("D9").Value = 0.02796733855954869
$ws.Range("E9").Value = 0.5331519221147261
$ws.Range("F9").Value = 0.6215739468373158
$ws.Range("I9").Value = 0.4458617095770414
$ws.Range("K9").Value = 0.4585399660351186
$ws.Range("N9").Value = 1.161075244190121
$ws.Range("O9").Value = 2.056200702562109
$ws.Range("B10").Value = 0.4879223112292266
$ws.Range("C10").Value = 0.09486110442108497
$ws.Range("D10").Value = 0.03147285446512882
$ws.Range("E10").Value = 0.6317487387890992
$ws.Range("F10").Value = 0.6313165380491554
$ws.Range("I10").Value = 0.4419975978579735
$ws.Range("K10").Value = 0.5413386479536939
$ws.Range("N10").Value = 1.139043985527552
$ws.Range("O10").Value = 2.064287731982006
$ws.Range("B11").Value = 0.520819553017958
$ws.Range("C11").Value = 0.1011895636160887
$ws.Range("D11").Value = 0.03305645241735533
$ws.Range("E11").Value = 0.6767467164284398
$ws.Range("F11").Value = 0.6361882923577724
$ws.Range("I11").Value = 0.4405267024191382
$ws.Range("K11").Value = 0.5789558663398964
$ws.Range("N11").Value = 1.129490597578853
$ws.Range("O11").Value = 2.069389979269175
$ws.Range("B12").Value = 0.5332744349278471
$ws.Range("C12").Value = 0.1035826445485384
$ws.Range("D12").Value = 0.03365449761152917
$ws.Range("E12").Value = 0.6938090645357846
$ws.Range("F12").Value = 0.6380964993384595
$ws.Range("I12").Value = 0.4400110636896635
$ws.Range("K12").Value = 0.5931932769146329
$ws.Range("N12").Value = 1.125940445706261
$ws.Range("O12").Value = 2.071527350385168
$ws.Range("B13").Value = 0.5305921762110586
$ws.Range("C13").Value = 0.103067402013977
$ws.Range("D13").Value = 0.03352577089765418
$ws.Range("E13").Value = 0.6901333518449064
$ws.Range("F13").Value = 0.6376827118321842
$ws.Range("I13").Value = 0.4401202747287201
$ws.Range("K13").Value = 0.590127332932127
$ws.Range("N13").Value = 1.126702030265719
$ws.Range("O13").Value = 2.071057891055972
$ws.Range("B14").Value = 0.5218442783456965
$ws.Range("C14").Value = 0.1013865114984753
$ws.Range("D14").Value = 0.03310568680723236
$ws.Range("E14").Value = 0.6781499839690497
$ws.Range("F14").Value = 0.6363440106950975
$ws.Range("I14").Value = 0.4404834510691416
$ws.Range("K14").Value = 0.5801273379196061
$ws.Range("N14").Value = 1.129197170862483
$ws.Range("O14").Value = 2.069561705041934
$ws.Range("B15").Value = 0.5164855865738787
$ws.Range("C15").Value = 0.1003564779753958
$ws.Range("D15").Value = 0.03284815996775592
$ws.Range("E15").Value = 0.6708128155605522
$ws.Range("F15").Value = 0.6355322755128583
$ws.Range("I15").Value = 0.4407112960116883
$ws.Range("K15").Value = 0.5740010729573441
$ws.Range("N15").Value = 1.130734311310619
$ws.Range("O15").Value = 2.068671996034084
$ws.Range("B16").Value = 0.4857720633426652
$ws.Range("C16").Value = 0.09444705468365555
$ws.Range("D16").Value = 0.03136913661393947
$ws.Range("E16").Value = 0.6288110875300958
$ws.Range("F16").Value = 0.6310070190307187
$ws.Range("I16").Value = 0.4420995062339905
$ws.Range("K16").Value = 0.5388792621643574
$ws.Range("N16").Value = 1.139677763028056
$ws.Range("O16").Value = 2.063982980021763
$ws.Range("B17").Value = 0.4669262220450321
$ws.Range("C17").Value = 0.09081584193900483
$ws.Range("D17").Value = 0.03045894059707877
$ws.Range("E17").Value = 0.6030828814736537
$ws.Range("F17").Value = 0.6283436613675519
$ws.Range("I17").Value = 0.4430246836350946
$ws.Range("K17").Value = 0.5173204619229637
$ws.Range("N17").Value = 1.145284410274492
$ws.Range("O17").Value = 2.0614714002987
$ws.Range("B18").Value = 0.4560852468399332
$ws.Range("C18").Value = 0.08872507784579398
$ws.Range("D18").Value = 0.029934379268866
$ws.Range("E18").Value = 0.588298322497721
$ws.Range("F18").Value = 0.6268531564255611
$ws.Range("I18").Value = 0.4435838203440596
$ws.Range("K18").Value = 0.5049158901428257
$ws.Range("N18").Value = 1.148553319260589
$ws.Range("O18").Value = 2.060160750249793
$ws.Range("B19").Value = 0.4524144626958275
$ws.Range("C19").Value = 0.08801680630992337
$ws.Range("D19").Value = 0.02975659417741383
$ws.Range("E19").Value = 0.5832948155520228
$ws.Range("F19").Value = 0.6263556013268712
$ws.Range("I19").Value = 0.4437777685057149
$ws.Range("K19").Value = 0.5007151531562783
$ws.Range("N19").Value = 1.149667689139234
$ws.Range("O19").Value = 2.059739974233082
$ws.Range("B20").Value = 0.4689325368490529
$ws.Range("C20").Value = 0.09120261722850387
$ws.Range("D20").Value = 0.03055594048053223
$ws.Range("E20").Value = 0.6058202711047045
$ws.Range("F20").Value = 0.6286228958502207
$ws.Range("I20").Value = 0.4429234018844532
$ws.Range("K20").Value = 0.5196159042021122
$ws.Range("N20").Value = 1.144683006743435
$ws.Range("O20").Value = 2.061724895360271
$ws.Range("B21").Value = 0.5244138221738694
$ws.Range("C21").Value = 0.1018803214363402
$ws.Range("D21").Value = 0.03322912018343516
$ws.Range("E21").Value = 0.6816691634750356
$ws.Range("F21").Value = 0.6367354984828353
$ws.Range("I21").Value = 0.4403756541427626
$ws.Range("K21").Value = 0.5830647837671847
$ws.Range("N21").Value = 1.128462454569926
$ws.Range("O21").Value = 2.069995595826327
$ws.Range("B22").Value = 0.5606587290529319
$ws.Range("C22").Value = 0.1088391680272309
$ws.Range("D22").Value = 0.03496668612160647
$ws.Range("E22").Value = 0.7313732378803905
$ws.Range("F22").Value = 0.6424070401718325
$ws.Range("I22").Value = 0.4389516634961055
$ws.Range("K22").Value = 0.6244888668702799
$ws.Range("N22").Value = 1.118254929640164
$ws.Range("O22").Value = 2.076597662456379
$ws.Range("B23").Value = 0.5413157032594427
$ws.Range("C23").Value = 0.1051269072553964
$ws.Range("D23").Value = 0.03404019661113722
$ws.Range("E23").Value = 0.7048325792939352
$ws.Range("F23").Value = 0.6393461783937369
$ws.Range("I23").Value = 0.4396895792600937
$ws.Range("K23").Value = 0.6023841812907165
$ws.Range("N23").Value = 1.123666833282712
$ws.Range("O23").Value = 2.072964327735207
$ws.Range("B24").Value = 0.4680255013583121
$ws.Range("C24").Value = 0.09102776587369021
$ws.Range("D24").Value = 0.0305120908092178
$ws.Range("E24").Value = 0.6045826758208079
$ws.Range("F24").Value = 0.6284965271862148
$ws.Range("I24").Value = 0.442969106518138
$ws.Range("K24").Value = 0.5185781662694922
$ws.Range("N24").Value = 1.144954759324653
$ws.Range("O24").Value = 2.061609875080507
$ws.Range("B25").Value = 0.3889203981696028
$ws.Range("C25").Value = 0.07573601648925887
$ws.Range("D25").Value = 0.02666680388574605
$ws.Range("E25").Value = 0.496960921705778
$ws.Range("F25").Value = 0.6122752868145653
$ws.Range("I25").Value = 0.447538429811452
$ws.Range("K25").Value = 0.4280085604366093
$ws.Range("N25").Value = 1.169603452473755
$ws.Range("O25").Value = 2.054480762449003
